# Update column G ("K") values for rows 2-36 on the active sheet.
# This regenerates the saved K (strikeouts-based) stat using the new
# calculation instead of the previous "Strike#" derived value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newKValues = @{
    2  = 4
    3  = 10
    4  = 9
    5  = 9
    6  = 16
    7  = 7
    8  = 9
    9  = 4
    10 = 5
    11 = 13
    12 = 6
    13 = 4
    14 = 8
    15 = 12
    16 = 11
    17 = 7
    18 = 4
    19 = 11
    20 = 8
    21 = 5
    22 = 7
    23 = 5
    24 = 0
    25 = 3
    26 = 3
    27 = 2
    28 = 8
    29 = 5
    30 = 0
    31 = 3
    32 = 3
    33 = 5
    34 = 5
    35 = 2
    36 = 3
}

foreach ($row in $newKValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $newKValues[$row]
}
